$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'68.259.18"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.43%  "
$ws.Range("D3").Value = "'3.818.70"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.52%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'600.19"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.32%  "
$ws.Range("D6").Value = "'169.25"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.47%  "
$ws.Range("D7").Value = "'3.820.43"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.46%  "
$ws.Range("E8").Value = "  -0.15%  "
$ws.Range("E9").Value = "  -1.86%  "
$ws.Range("E10").Value = "  -3.16%  "
$ws.Range("D11").Value = "'6.46"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.83%  "
$ws.Range("E12").Value = "  -2.38%  "
$ws.Range("D13").Value = "'0.0000261"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.83%  "
$ws.Range("D14").Value = "'37.02"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.55%  "
$ws.Range("D15").Value = "'4.458.31"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.54%  "
$ws.Range("D16").Value = "'3.813.68"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.72%  "
$ws.Range("D17").Value = "'68.266.30"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.43%  "
$ws.Range("D18").Value = "'18.55"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.84%  "
$ws.Range("E19").Value = "  -2.88%  "
$ws.Range("D21").Value = "'11.12"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "'468.13"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.33%  "
$ws.Range("D23").Value = "'0.733"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.85%  "
$ws.Range("D24").Value = "'0.0000160"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.25%  "
$ws.Range("D25").Value = "'82.98"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.20%  "
$ws.Range("E26").Value = "  -3.12%  "
$ws.Range("D27").Value = "'12.12"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.63%  "
$ws.Range("D28").Value = "'10.01"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.38%  "
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("E30").Value = "  -1.28%  "
$ws.Range("D31").Value = "'3.966.77"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.48%  "
$ws.Range("D32").Value = "'7.67"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.13%  "
$ws.Range("D33").Value = "'31.44"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.31%  "
$ws.Range("E34").Value = "  -5.17%  "
$ws.Range("D35").Value = "'9.45"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.16%  "
$ws.Range("D36").Value = "'3.781.08"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.56%  "
$ws.Range("E37").Value = "  -3.30%  "
$ws.Range("E38").Value = "  +11.73%  "
$ws.Range("E39").Value = "  -1.28%  "
$ws.Range("E40").Value = "  -2.93%  "
$ws.Range("E41").Value = "  -3.45%  "
$ws.Range("D42").Value = "'1.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.09%  "
$ws.Range("D43").Value = "'0.314"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.68%  "
$ws.Range("E44").Value = "  -6.45%  "
$ws.Range("E45").Value = "  +0.73%  "
$ws.Range("D46").Value = "'0.000295"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +7.51%  "
$ws.Range("D47").Value = "'417.36"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.55%  "
$ws.Range("E48").Value = "  -0.02%  "
$ws.Range("D49").Value = "'46.97"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.64%  "
$ws.Range("D50").Value = "'26.19"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.25%  "
$ws.Range("D51").Value = "'141.75"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.96%  "
